$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stimulus")

# Row 19 (Test No 16) currently holds the "Edge Cases:Burst with maximum
# transfer size (HSIZE = WORD)" entry. Remove it entirely: rows below
# shift up one, and the now-unused shared string is dropped automatically.
$ws.Rows(19).Delete()

# Renumber the "Test No" column for the rows that shifted up so the
# sequence stays contiguous (1..17).
$ws.Range("B19").Value = 16
$ws.Range("B20").Value = 17

# Update the view selection to match the saved state.
$ws.Range("B18:D18").Select() | Out-Null
